# Update column F (dSF) values for specific rows based on repulled/recomputed
# data (commit message: "repull data, push all data, mean calculation").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 3
    17 = 3
    25 = -3
    28 = 3
    41 = 0
    43 = -1
    44 = -1
    46 = 0
    47 = -8
    49 = -2
    51 = 9
    54 = -4
    57 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
